$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the row for ID "RM 232" (originally row 26).
$ws.Rows.Item(26).Delete()

# After that deletion, the row for ID "SC 92" (originally row 28) is now row 27.
$ws.Rows.Item(27).Delete()

# Update column F (and a few column C) values that moved between
# "missing" (blank) and "present" (numeric) between the two runs.
$ws.Range("F2").Value = 18.03
$ws.Range("F12").Value = 17.45
$ws.Range("F20").Value = 17.73
$ws.Range("F21").Value = 16.58
$ws.Range("C26").Value = 10.8
$ws.Range("C29").Value = 11.2
$ws.Range("C30").Value = 11.4
$ws.Range("F31").Value = 17.18
$ws.Range("F33").Value = 17.53

$ws.Range("F6").ClearContents()
$ws.Range("F14").ClearContents()
$ws.Range("F23").ClearContents()
$ws.Range("F24").ClearContents()
$ws.Range("C27").ClearContents()
$ws.Range("C28").ClearContents()
$ws.Range("C31").ClearContents()
$ws.Range("C32").ClearContents()
